$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BILLING INVOICE")

# H9: invoice/billing date moves forward 5 days
$ws.Range("H9").Value = 45289.66666666667

# Row 16 (existing internal salary/deduction placeholder row)
$ws.Range("B16").Value = 45265.298125
$ws.Range("H16").Value = 500

# Row 18: new placeholder row - sequence number, date, shared-string label, amount
$ws.Range("A18").Value = 2
$ws.Range("B18").Value = 45265.29819444445
$ws.Range("D18").Value = "SHIPMENT NO.: 123123234124"
$ws.Range("H18").Value = 5555

# Row 19: new placeholder label only
$ws.Range("D19").Value = "SPO NO.: 1241421414"

# Writing into the accounting-formatted H column cells re-triggers this
# sheet's row autofit; restore the original (manually set) row heights so
# only the intended cell values change.
$ws.Rows(16).RowHeight = 9.75
$ws.Rows(18).RowHeight = 9.75
